$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (B2:I2) to all have the same accuracy value
$ws.Range("B2:I2").Value = 87.36702129244804
